$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Add a new record (row 3) to the table on Hoja1
$ws.Range("A3").Value = "QF13002"
$ws.Range("B3").Value = "Aristides"
$ws.Range("C3").Value = "Fuentes"
$ws.Range("D3").Value = 28
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 7.5
$ws.Range("G3").Value = 2013
$ws.Range("H3").Value = 7
$ws.Range("I3").Value = 1

# Match the selection that results from entering data and landing on I3
$ws.Range("I3").Select()
